$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 35 with date and error count, matching style of row 34 (A column date-formatted)
$ws.Range("A34:B34").Copy()
$ws.Range("A35:B35").PasteSpecial(-4122)

$ws.Range("A35").Value = 46007
$ws.Range("B35").Value = 50

# Update selection to mirror the new last row
$ws.Range("A35:B35").Select()
